$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting (avoid numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values from the crypto data refresh
$ws.Range("D2").Value = "63.017.41"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "3.471.22"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "578.08"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "148.44"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("D7").Value = "3.471.65"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.480"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "0.403"
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("D13").Value = "4.064.10"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "29.78"
$ws.Range("E14").Value = "  +6.16%  "
$ws.Range("D16").Value = "3.471.08"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "62.953.63"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("D20").Value = "14.39"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("D21").Value = "9.23"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "388.51"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "0.558"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "74.67"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "3.609.42"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "8.15"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "23.69"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D36").Value = "32.17"
$ws.Range("E36").Value = "  +18.73%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "5.28"
$ws.Range("E37").Value = "  +3.49%  "
$ws.Range("D38").Value = "7.08"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "169.69"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  +5.41%  "
$ws.Range("D41").Value = "3.507.84"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").Value = "0.0756"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").Value = "0.799"
$ws.Range("D44").Value = "42.39"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "4.46"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").Value = "1.72"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  +4.28%  "
$ws.Range("D48").Value = "2.619.24"
$ws.Range("E48").Value = "  +5.60%  "
$ws.Range("E49").Value = "  +11.92%  "
$ws.Range("D50").Value = "22.99"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").Value = "6.74"
$ws.Range("E51").Value = "  +1.10%  "
